$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Sheet is protected; temporarily unprotect to allow the edits, then
# restore protection with the same password afterwards.
$ws.Unprotect("CC21")

# "Serviço atuais:" -> "Serviço:"
$ws.Range("B3").Value = "Serviço:"

# Remove the old "tela cidade" entry from the "Telas já prontas" list
$ws.Range("E7").Value = ""

# Add the new "tela cad. Cidade" service entry
$ws.Range("B11").Value = "tela cad. Cidade"

$ws.Protect("CC21")

# Update selection to match the author's final cursor position
$ws.Range("B18").Select()
